$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.016.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").Value = "'1.554.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "'290.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").Value = "'0.3968"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.84%  "

$ws.Range("D8").Value = "'0.3224"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.11%  "

$ws.Range("D9").Value = "'43.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("D10").Value = "'0.07255"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("D11").Value = "'1.079"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.32%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "'5.706"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("D14").Value = "'18.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.99%  "

$ws.Range("D15").Value = "'0.00001129"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.13%  "

$ws.Range("D16").Value = "'6.630"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "'1.555.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").Value = "'0.06578"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.15%  "

$ws.Range("D19").Value = "'83.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.37%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'6.273"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "

$ws.Range("D22").Value = "'15.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.74%  "

$ws.Range("D23").Value = "'11.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("D24").Value = "'22.030.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Value = "'2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.11%  "

$ws.Range("D26").Value = "'2.415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.38%  "

$ws.Range("D27").Value = "'148.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("D28").Value = "'18.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.63%  "

$ws.Range("D29").Value = "'4.871"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").Value = "'1.726.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").Value = "'118.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("D32").Value = "'0.9674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.31%  "

$ws.Range("D33").Value = "'5.826"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").Value = "'0.08315"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").Value = "'9.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").Value = "'1.599"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -16.73%  "

$ws.Range("D37").Value = "'0.02262"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.62%  "

$ws.Range("D38").Value = "'5.108"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("D39").Value = "'0.05992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.19%  "

$ws.Range("D40").Value = "'1.215"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").Value = "'0.2037"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.64%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "'10.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "

$ws.Range("D44").Value = "'0.5814"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.13%  "

$ws.Range("D45").Value = "'13.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.92%  "

$ws.Range("D46").Value = "'3.744"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "'0.5574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.96%  "

$ws.Range("D48").Value = "'118.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("D49").Value = "'1.899"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("D50").Value = "'1.133"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").Value = "'0.06814"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.46%  "
